$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.175.83"
$ws.Range("E2").Value = "  +0.98%  "
$ws.Range("D3").Value = "3.855.79"
$ws.Range("E3").Value = "  +1.07%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "696.39"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.83%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.21"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.01%  "
$ws.Range("D7").Value = "3.853.72"
$ws.Range("E7").Value = "  +1.06%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  +0.31%  "
$ws.Range("E10").Value = "  +2.00%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.35"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.06%  "
$ws.Range("E12").Value = "  +0.40%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000259"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.62"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.17%  "
$ws.Range("D15").Value = "4.503.26"
$ws.Range("E15").Value = "  +1.01%  "
$ws.Range("D16").Value = "3.864.64"
$ws.Range("E16").Value = "  +1.34%  "
$ws.Range("D17").Value = "71.178.92"
$ws.Range("E17").Value = "  +1.00%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.77"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.73%  "
$ws.Range("E19").Value = "  +0.98%  "
$ws.Range("E20").Value = "  +0.36%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.16"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.27%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "496.35"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.22%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.725"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.65%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.03"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.92%  "
$ws.Range("E25").Value = "  +4.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.36"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.60"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.08%  "
$ws.Range("E28").Value = "  +2.40%  "
$ws.Range("D29").Value = "4.010.16"
$ws.Range("E29").Value = "  +1.12%  "
$ws.Range("E30").Value = "  +10.74%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.64"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.44%  "
$ws.Range("E33").Value = "  +0.43%  "
$ws.Range("E34").Value = "  +0.76%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.179"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.30%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.32"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.48%  "
$ws.Range("D37").Value = "3.807.55"
$ws.Range("E37").Value = "  +0.95%  "
$ws.Range("E38").Value = "  -0.13%  "
$ws.Range("E39").Value = "  +3.44%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.39"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +13.38%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.07"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.02%  "
$ws.Range("E43").Value = "  +5.82%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.999"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.11%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "164.41"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.37%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.000307"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.75%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "48.69"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.20%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "44.32"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.22%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "420.26"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.89%  "
$ws.Range("E51").Value = "  +0.90%  "
